# Generate Report for Handoff
# The CI run replaced the two tracked e2e markdown files with a new pair
# (new GUID file names) and refreshed their localization-status rows:
#   - file/path cells now point at the new source files
#   - status moved from "Handed back: in sync with en-US" to "Ready for handoff"
#   - handoff timestamps/hashes refreshed, handback (not yet done) cleared
#   - zh-cn sheet's second row duplicate flag flipped True, one xlf re-hashed
# This script reproduces that row-level refresh plus the resulting column
# width changes across all three sheets.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "4a1a6148-7426-402f-b4cd-6613294a232b"
$oldGuid2 = "ee03db45-3948-4fe5-a1cb-902e4795ba7e"
$newGuid1 = "76c58e28-84cf-41ee-b692-2bb5c64269f6"
$newGuid2 = "ffffea95f99f-e464-439d-805b-ad2129f8bb03"

$newFile1 = "$newGuid1.md"
$newFile2 = "$newGuid2.md"
$newPath1 = "e2e\$newFile1"
$newPath2 = "e2e\$newFile2"

$urlBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6ab939b926161d86e4d338e012f20a3db3d60509/e2e/"
$url1 = "$urlBase$newFile1"
$url2 = "$urlBase$newFile2"

$status = "Ready for handoff"
$hoDate = "2016-08-29 01:01:09"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value2 = $newFile1
$wsOverview.Range("E2").Value2 = $status
$wsOverview.Range("F2").Value2 = $status
$wsOverview.Range("G2").Value2 = $hoDate

$wsOverview.Range("A3").Value2 = $newFile2
$wsOverview.Range("E3").Value2 = $status
$wsOverview.Range("F3").Value2 = $status
$wsOverview.Range("G3").Value2 = $hoDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $url1, [System.Type]::Missing, [System.Type]::Missing, $newPath1) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $url2, [System.Type]::Missing, [System.Type]::Missing, $newPath2) | Out-Null

$wsOverview.Columns.Item(5).ColumnWidth = 16.333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhHandoff1 = "$newGuid1.a650dce89da627fa1327405f9da5a0e0a35480f4.zh-cn.xlf"
$zhHandoffDate = "2016-08-29 01:01:01"
$zeroDate = "0001-01-01 00:00:00"

$wsZh.Range("C2").Value2 = $status
$wsZh.Range("G2").Value2 = $zhHandoff1
$wsZh.Range("H2").Value2 = $zhHandoffDate
$wsZh.Range("I2").Value2 = ""
$wsZh.Range("J2").Value2 = ""
$wsZh.Range("K2").Value2 = $zeroDate

$wsZh.Range("C3").Value2 = $status
$wsZh.Range("F3").Value2 = "True"
$wsZh.Range("G3").Value2 = $zhHandoff1
$wsZh.Range("H3").Value2 = $zhHandoffDate
$wsZh.Range("I3").Value2 = ""
$wsZh.Range("J3").Value2 = ""
$wsZh.Range("K3").Value2 = $zeroDate

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $url1, [System.Type]::Missing, [System.Type]::Missing, $newFile1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $url2, [System.Type]::Missing, [System.Type]::Missing, $newFile2) | Out-Null

$wsZh.Columns.Item(3).ColumnWidth = 16.333333
$wsZh.Columns.Item(9).ColumnWidth = 17.833333
$wsZh.Columns.Item(10).ColumnWidth = 20.833333

# reset I/J columns back to the default (non-hyperlink) style now that
# their hyperlinks are gone
$wsZh.Range("I2:J3").Font.Underline = $false
$wsZh.Range("I2:J3").Font.ColorIndex = 1

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deHandoff1 = "$newGuid1.a650dce89da627fa1327405f9da5a0e0a35480f4.de-de.xlf"

$wsDe.Range("C2").Value2 = $status
$wsDe.Range("G2").Value2 = $deHandoff1
$wsDe.Range("H2").Value2 = $hoDate
$wsDe.Range("I2").Value2 = ""
$wsDe.Range("J2").Value2 = ""
$wsDe.Range("K2").Value2 = $zeroDate

$wsDe.Range("C3").Value2 = $status
$wsDe.Range("F3").Value2 = "True"
$wsDe.Range("G3").Value2 = $deHandoff1
$wsDe.Range("H3").Value2 = $hoDate
$wsDe.Range("I3").Value2 = ""
$wsDe.Range("J3").Value2 = ""
$wsDe.Range("K3").Value2 = $zeroDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $url1, [System.Type]::Missing, [System.Type]::Missing, $newFile1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $url2, [System.Type]::Missing, [System.Type]::Missing, $newFile2) | Out-Null

$wsDe.Columns.Item(3).ColumnWidth = 16.333333
$wsDe.Columns.Item(9).ColumnWidth = 17.833333
$wsDe.Columns.Item(10).ColumnWidth = 20.833333

$wsDe.Range("I2:J3").Font.Underline = $false
$wsDe.Range("I2:J3").Font.ColorIndex = 1
